# Atualização de bases das ligas, do dia: 17-02-2024 às 11:11
#
# The source data rows got re-sorted within same-date match blocks; for a
# handful of dates this changed which physical row holds which match's
# data. Column A (running index) stays put per-row; columns B and F:AC
# (the actual match data) rotate among the rows of each affected group -
# each row ends up holding the data that previously lived in the row
# immediately above it (wrapping around within the group).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row groups (1-based worksheet rows) that are cyclically rotated.
$groups = @(
    @(14, 15),
    @(25, 26, 27),
    @(39, 40),
    @(50, 51),
    @(81, 82, 83),
    @(108, 109),
    @(124, 125),
    @(139, 140),
    @(155, 156),
    @(176, 177),
    @(204, 205, 206, 207)
)

foreach ($rows in $groups) {
    $n = $rows.Count

    # Capture the "B" and "F:AC" data for every row in the group before
    # writing anything (so later writes don't clobber data we still need
    # to read).
    $bVals = @()
    $facVals = @()
    foreach ($r in $rows) {
        $bVals += , ($ws.Range("B$r").Value2)
        $facVals += , ($ws.Range("F" + $r + ":AC" + $r).Value2)
    }

    # Each row takes on the data previously held by the row above it in
    # the group (wrapping around).
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $rows[$i]
        $srcIdx = ($i - 1 + $n) % $n

        $ws.Range("B$destRow").Value2 = $bVals[$srcIdx]
        $ws.Range("F" + $destRow + ":AC" + $destRow).Value2 = $facVals[$srcIdx]
    }
}
